# Update PLC data 2025-10-13 14:02:07
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 172935
$ws.Range("C4").Value = 163719
$ws.Range("C5").Value = 9216
$ws.Range("C6").Value = 755
$ws.Range("C8").Value = 66.15000000000001
